$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "waiter"
$ws.Range("B2").Value = "ウエイター"
$ws.Range("A3").Value = "(someone's) house/home"
$ws.Range("B3").Value = "お宅|おたく"
$ws.Range("A4").Value = "adult"
$ws.Range("B4").Value = "大人|おとな"
$ws.Range("A5").Value = "foreign language"
$ws.Range("B5").Value = "外国語|がいこくご"
$ws.Range("A6").Value = "musical instrument"
$ws.Range("B6").Value = "楽器|がっき"
$ws.Range("A7").Value = "karate"
$ws.Range("B7").Value = "空手|からて"
$ws.Range("A8").Value = "curry"
$ws.Range("B8").Value = "カレー"
$ws.Range("A9").Value = "kimono; Japanese traditional dress"
$ws.Range("B9").Value = "着物|きもの"
$ws.Range("A10").Value = "advertisement"
$ws.Range("B10").Value = "広告|こうこく"
$ws.Range("A11").Value = "tea (black tea)"
$ws.Range("B11").Value = "紅茶|こうちゃ"
$ws.Range("A12").Value = "language"
$ws.Range("B12").Value = "言葉|ことば"
$ws.Range("A13").Value = "golf"
$ws.Range("B13").Value = "ゴルフ"
$ws.Range("A14").Value = "sweater"
$ws.Range("B14").Value = "セーター"
$ws.Range("A15").Value = "elephant"
$ws.Range("B15").Value = "象|ぞう"
$ws.Range("A16").Value = "violin"
$ws.Range("B16").Value = "バイオリン"
$ws.Range("A17").Value = "motorcycle"
$ws.Range("B17").Value = "バイク"
$ws.Range("A18").Value = "(consumer) prices"
$ws.Range("B18").Value = "物価|ぶっか"
$ws.Range("A19").Value = "grammar"
$ws.Range("B19").Value = "文法|ぶんぽう"
$ws.Range("A20").Value = "lawyer"
$ws.Range("B20").Value = "弁護士|べんごし"
$ws.Range("A21").Value = "recruitment"
$ws.Range("B21").Value = "募集|ぼしゅう"
$ws.Range("A22").Value = "shop; store"
$ws.Range("B22").Value = "店|みせ"
$ws.Range("A23").Value = "yakuza; gangster"
$ws.Range("B23").Value = "やくざ"
$ws.Range("A24").Value = "promise; appointment"
$ws.Range("B24").Value = "約束|やくそく"
$ws.Range("A25").Value = "(term) paper"
$ws.Range("B25").Value = "レポート"
$ws.Range("A26").Value = "I (formal)"
$ws.Range("B26").Value = "私|わたくし"
$ws.Range("A27").Value = "glad"
$ws.Range("B27").Value = "うれしい"
$ws.Range("A28").Value = "sad"
$ws.Range("B28").Value = "悲しい|かなしい"
$ws.Range("A29").Value = "hot and spicy; salty"
$ws.Range("B29").Value = "辛い|からい"
$ws.Range("A30").Value = "strict"
$ws.Range("B30").Value = "厳しい|きびしい"
$ws.Range("A31").Value = "incredible; awesome"
$ws.Range("B31").Value = "すごい"
$ws.Range("A32").Value = "close; near"
$ws.Range("B32").Value = "近い|ちかい"
$ws.Range("A33").Value = "various; different kinds of"
$ws.Range("B33").Value = "いろいろ（な）"
$ws.Range("A34").Value = "happy (lasting happiness)"
$ws.Range("B34").Value = "幸せ|しあわせ（な）"
$ws.Range("A35").Value = "no good"
$ws.Range("B35").Value = "だめ（な）"
$ws.Range("A36").Value = "to knit"
$ws.Range("B36").Value = "編む|あむ"
$ws.Range("A37").Value = "to lend; to rent"
$ws.Range("B37").Value = "貸す|かす"
$ws.Range("A38").Value = "to do one's best; to try hard"
$ws.Range("B38").Value = "頑張る|がんばる"
$ws.Range("A39").Value = "to cry"
$ws.Range("B39").Value = "泣く|なく"
$ws.Range("A40").Value = "to brush (teeth); to polish"
$ws.Range("B40").Value = "磨く|みがく"
$ws.Range("A41").Value = "to keep a promise"
$ws.Range("B41").Value = "約束を守る|やくそくをまもる"
$ws.Range("A42").Value = "to be moved/touched (by...)"
$ws.Range("B42").Value = "感動する|かんどうする"
$ws.Range("A43").Value = "someone honorable is present/home"
$ws.Range("B43").Value = "いらっしゃいませ"
$ws.Range("A44").Value = "...times"
$ws.Range("B44").Value = "～回|～かい"
$ws.Range("A45").Value = "...kilometers; ...kilograms"
$ws.Range("B45").Value = "～キロ"
$ws.Range("A46").Value = "as a matter of fact,..."
$ws.Range("B46").Value = "実は|じつは"
$ws.Range("A47").Value = "all"
$ws.Range("B47").Value = "全部|ぜんぶ"
$ws.Range("A48").Value = "my name is..."
$ws.Range("B48").Value = "～と申します|～ともうします"
$ws.Range("A49").Value = "one day"
$ws.Range("B49").Value = "一日|いちにち"
$ws.Range("A50").Value = "two days"
$ws.Range("B50").Value = "二日|ふつか"
$ws.Range("A51").Value = "three days"
$ws.Range("B51").Value = "三日|みっか"
$ws.Range("A52").Value = "four days"
$ws.Range("B52").Value = "四日|よっか"
$ws.Range("A53").Value = "five days"
$ws.Range("B53").Value = "五日|いつか"
$ws.Range("A54").Value = "six days"
$ws.Range("B54").Value = "六日|むいか"
$ws.Range("A55").Value = "seven days"
$ws.Range("B55").Value = "七日|なのか"
$ws.Range("A56").Value = "eight days"
$ws.Range("B56").Value = "八日|ようか"
$ws.Range("A57").Value = "nine days"
$ws.Range("B57").Value = "九日|ここのか"
$ws.Range("A58").Value = "ten days"
$ws.Range("B58").Value = "十日|とおか"
$ws.Range("A73").Value = "I would like to open an account."
$ws.Range("B73").Value = "口座を開きたいんですが。|こうざをひらきたいんですが。"
$ws.Range("A74").Value = "I would like to close an account."
$ws.Range("B74").Value = "口座を閉じたいんですが。|こうざをとじたいんですが。"
$ws.Range("A75").Value = "Please change dollars into yen."
$ws.Range("B75").Value = "ドルを円にかえてください。|ドルを円にかえてください。"
$ws.Range("A76").Value = "I would like to deposit money into the account."
$ws.Range("B76").Value = "口座にお金を振り込みたいんですが。|こうざにおかねをふりこみたいんですが。"
$ws.Range("A77").Value = "Can you change a 10,000-yen bill into ten 1,000-yen bills?"
$ws.Range("B77").Value = "一万円札を千円札十枚に両替できますか。|いちまんえんさつをせんえんさつじゅうまいにりょうがえできますか。"
$ws.Range("A78").Value = "I will withdraw money."
$ws.Range("B78").Value = "お金をおろします。|おかねをおろします。"
$ws.Range("A79").Value = "food"
$ws.Range("B79").Value = "食べ物|たべもの"
$ws.Range("A80").Value = "drink"
$ws.Range("B80").Value = "飲み物|のみもの"
$ws.Range("A81").Value = "things"
$ws.Range("B81").Value = "物|もの"
$ws.Range("A82").Value = "shopping"
$ws.Range("B82").Value = "買い物|かいもの"
$ws.Range("A83").Value = "animal"
$ws.Range("B83").Value = "動物|どうぶつ"
$ws.Range("A84").Value = "bird"
$ws.Range("B84").Value = "鳥|とり"
$ws.Range("A85").Value = "grilled chicken"
$ws.Range("B85").Value = "焼き鳥|やきとり"
$ws.Range("A86").Value = "swan"
$ws.Range("B86").Value = "白鳥|はくちょう"
$ws.Range("A87").Value = "cooking"
$ws.Range("B87").Value = "料理|りょうり"
$ws.Range("A88").Value = "charge"
$ws.Range("B88").Value = "料金|りょうきん"
$ws.Range("A89").Value = "tuition"
$ws.Range("B89").Value = "授業料|じゅぎょうりょう"
$ws.Range("A90").Value = "salary"
$ws.Range("B90").Value = "給料|きゅうりょう"
$ws.Range("A91").Value = "cooking"
$ws.Range("B91").Value = "料理|りょうり"
$ws.Range("A92").Value = "reason"
$ws.Range("B92").Value = "理由|りゆう"
$ws.Range("A93").Value = "geography"
$ws.Range("B93").Value = "地理|ちり"
$ws.Range("A94").Value = "impossible"
$ws.Range("B94").Value = "無理な|むりな"
$ws.Range("A95").Value = "especially"
$ws.Range("B95").Value = "特に|とくに"
$ws.Range("A96").Value = "special"
$ws.Range("B96").Value = "特別な|とくべつな"
$ws.Range("A97").Value = "characteristic"
$ws.Range("B97").Value = "特徴|とくちょう"
$ws.Range("A98").Value = "super express"
$ws.Range("B98").Value = "特急|とっきゅう"
$ws.Range("A99").Value = "cheap"
$ws.Range("B99").Value = "安い|やすい"
$ws.Range("A100").Value = "safe"
$ws.Range("B100").Value = "安全な|あんぜんな"
$ws.Range("A101").Value = "relief"
$ws.Range("B101").Value = "安心|あんしん"
$ws.Range("A102").Value = "uneasy"
$ws.Range("B102").Value = "不安な|ふあんな"
$ws.Range("A103").Value = "rice; meal"
$ws.Range("B103").Value = "ご飯|ごはん"
$ws.Range("A104").Value = "breakfast"
$ws.Range("B104").Value = "朝ご飯|あさごはん"
$ws.Range("A105").Value = "dinner"
$ws.Range("B105").Value = "晩ご飯|ばんごはん"
$ws.Range("A106").Value = "meat"
$ws.Range("B106").Value = "肉|にく"
$ws.Range("A107").Value = "beef"
$ws.Range("B107").Value = "牛肉|ぎゅうにく"
$ws.Range("A108").Value = "pork"
$ws.Range("B108").Value = "豚肉|ぶたにく"
$ws.Range("A109").Value = "meat shop"
$ws.Range("B109").Value = "肉屋|にくや"
$ws.Range("A110").Value = "muscle"
$ws.Range("B110").Value = "筋肉|きんにく"
$ws.Range("A111").Value = "bad"
$ws.Range("B111").Value = "悪い|わるい"
$ws.Range("A112").Value = "to feel sick"
$ws.Range("B112").Value = "気分が悪い|きぶんがわるい"
$ws.Range("A113").Value = "the worst"
$ws.Range("B113").Value = "最悪|さいあく"
$ws.Range("A114").Value = "devil"
$ws.Range("B114").Value = "悪魔|あくま"
$ws.Range("A115").Value = "body"
$ws.Range("B115").Value = "体|からだ"
$ws.Range("A116").Value = "body weight"
$ws.Range("B116").Value = "体重|たいじゅう"
$ws.Range("A117").Value = "gymnastics; physical exercises"
$ws.Range("B117").Value = "体操|たいそう"
$ws.Range("A118").Value = "airport"
$ws.Range("B118").Value = "空港|くうこう"
$ws.Range("A119").Value = "air"
$ws.Range("B119").Value = "空気|くうき"
$ws.Range("A120").Value = "sky"
$ws.Range("B120").Value = "空|そら"
$ws.Range("A121").Value = "to be vacant"
$ws.Range("B121").Value = "空く|あく"
$ws.Range("A122").Value = "karate"
$ws.Range("B122").Value = "空手|からて"
$ws.Range("A123").Value = "Kobe Port"
$ws.Range("B123").Value = "神戸港|こうべこう"
$ws.Range("A124").Value = "port"
$ws.Range("B124").Value = "港|みなと"
$ws.Range("A125").Value = "Hong Kong"
$ws.Range("B125").Value = "香港|ほんこん"
$ws.Range("A126").Value = "to arrive"
$ws.Range("B126").Value = "着く|つく"
$ws.Range("A127").Value = "to wear"
$ws.Range("B127").Value = "着る|きる"
$ws.Range("A128").Value = "kimono"
$ws.Range("B128").Value = "着物|きもの"
$ws.Range("A129").Value = "arriving at Osaka"
$ws.Range("B129").Value = "大阪着|おおさかちゃく"
$ws.Range("A130").Value = "the same"
$ws.Range("B130").Value = "同じ|おなじ"
$ws.Range("A131").Value = "coworker"
$ws.Range("B131").Value = "同僚|どうりょう"
$ws.Range("A132").Value = "classmate"
$ws.Range("B132").Value = "同級生|どうきゅうせい"
$ws.Range("A133").Value = "same time"
$ws.Range("B133").Value = "同時|どうじ"
$ws.Range("A134").Value = "sea"
$ws.Range("B134").Value = "海|うみ"
$ws.Range("A135").Value = "the Japan Sea"
$ws.Range("B135").Value = "日本海|にほんかい"
$ws.Range("A136").Value = "overseas"
$ws.Range("B136").Value = "海外|かいがい"
$ws.Range("A137").Value = "coast"
$ws.Range("B137").Value = "海岸|かいがん"
$ws.Range("A138").Value = "noon"
$ws.Range("B138").Value = "昼|ひる"
$ws.Range("A139").Value = "lunch"
$ws.Range("B139").Value = "昼ご飯|ひるごはん"
$ws.Range("A140").Value = "nap"
$ws.Range("B140").Value = "昼寝|ひるね"
$ws.Range("A141").Value = "lunch break"
$ws.Range("B141").Value = "昼休み|ひるやすみ"
$ws.Range("A142").Value = "lunch (formal)"
$ws.Range("B142").Value = "昼食|ちゅうしょく"
